# Generate Report for Handback
# Re-generates the handback-status workbook rows for the new file pair:
#   old "9b2e30a5-86ac-462a-ba19-7eb1ba5d53f2.md" -> new "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md"
#   old "bafaaef2-f42e-4bc8-a5c9-d3ed2d7fab02.md" -> new "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md"
# along with refreshed timestamps and xliff hash filenames.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md"
$overview.Range("B2").Value = "e2e\53392d28-b9dc-4a8c-adbd-6fd2226efb47.md"
$overview.Range("G2").Value = "2016-08-19 13:02:56"

$overview.Range("A3").Value = "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md"
$overview.Range("B3").Value = "e2e\ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md"
$overview.Range("G3").Value = "2016-08-19 13:02:56"

# Rebuild hyperlinks with unchanged target addresses but refreshed display text.
$overview.Range("A1").Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58015d1567b9974358bac4fd2f06aa791ea0ebb6/e2e/9b2e30a5-86ac-462a-ba19-7eb1ba5d53f2.md", "", "", "e2e\53392d28-b9dc-4a8c-adbd-6fd2226efb47.md")
$overview.Hyperlinks.Add($overview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58015d1567b9974358bac4fd2f06aa791ea0ebb6/e2e/bafaaef2-f42e-4bc8-a5c9-d3ed2d7fab02.md", "", "", "e2e\ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md"
$zhcn.Range("I2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md"
$zhcn.Range("G2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-19 13:02:51"
$zhcn.Range("J2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 13:03:16"

$zhcn.Range("A3").Value = "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md"
$zhcn.Range("I3").Value = "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md"
$zhcn.Range("G3").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-19 13:02:51"
$zhcn.Range("J3").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-19 13:03:16"

$zhcn.Range("A1").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58015d1567b9974358bac4fd2f06aa791ea0ebb6/e2e/9b2e30a5-86ac-462a-ba19-7eb1ba5d53f2.md", "", "", "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e990a7579950bcc7582cdb99f72670e28cefffdd/e2e/9b2e30a5-86ac-462a-ba19-7eb1ba5d53f2.md", "", "", "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58015d1567b9974358bac4fd2f06aa791ea0ebb6/e2e/bafaaef2-f42e-4bc8-a5c9-d3ed2d7fab02.md", "", "", "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e990a7579950bcc7582cdb99f72670e28cefffdd/e2e/bafaaef2-f42e-4bc8-a5c9-d3ed2d7fab02.md", "", "", "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md"
$dede.Range("I2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md"
$dede.Range("G2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.de-de.xlf"
$dede.Range("H2").Value = "2016-08-19 13:02:56"
$dede.Range("J2").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 13:03:23"

$dede.Range("A3").Value = "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md"
$dede.Range("I3").Value = "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md"
$dede.Range("G3").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.de-de.xlf"
$dede.Range("H3").Value = "2016-08-19 13:02:56"
$dede.Range("J3").Value = "53392d28-b9dc-4a8c-adbd-6fd2226efb47.1379c71ff5ec595ad43afa7e4fd9ef5a64c3a2f9.de-de.xlf"
$dede.Range("K3").Value = "2016-08-19 13:03:23"

$dede.Range("A1").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58015d1567b9974358bac4fd2f06aa791ea0ebb6/e2e/9b2e30a5-86ac-462a-ba19-7eb1ba5d53f2.md", "", "", "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3bc972a60344989a58d304a2b43c5ad678d3707e/e2e/9b2e30a5-86ac-462a-ba19-7eb1ba5d53f2.md", "", "", "53392d28-b9dc-4a8c-adbd-6fd2226efb47.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58015d1567b9974358bac4fd2f06aa791ea0ebb6/e2e/bafaaef2-f42e-4bc8-a5c9-d3ed2d7fab02.md", "", "", "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3bc972a60344989a58d304a2b43c5ad678d3707e/e2e/bafaaef2-f42e-4bc8-a5c9-d3ed2d7fab02.md", "", "", "ffff9b15097a-0a74-453e-8e14-f068e864bdcf.md")

$wb.Save()
